$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as plain text (e.g. "63.367.26", "1.00")
# rather than numbers, so force each target cell to keep a Text format before
# writing the new value -- otherwise Excel auto-coerces numeric-looking strings
# into actual numbers and drops formatting such as trailing zeros.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "63.367.26"
$ws.Cells.Item(2, 5).Value = "  -1.25%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.098.95"
$ws.Cells.Item(3, 5).Value = "  +0.93%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.06%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "556.52"
$ws.Cells.Item(5, 5).Value = "  +0.88%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "138.12"
$ws.Cells.Item(6, 5).Value = "  -2.43%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.04%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "3.092.85"
$ws.Cells.Item(8, 5).Value = "  +0.95%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.497"
$ws.Cells.Item(9, 5).Value = "  +1.52%  "

# Row 10
$ws.Cells.Item(10, 2).Value = "Toncoin"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(10, 4).Value = "6.72"
$ws.Cells.Item(10, 5).Value = "  +2.24%  "

# Row 11
$ws.Cells.Item(11, 2).Value = "Dogecoin"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(11, 4).Value = "0.162"
$ws.Cells.Item(11, 5).Value = "  +6.07%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "0.456"
$ws.Cells.Item(12, 5).Value = "  +0.61%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "35.14"
$ws.Cells.Item(13, 5).Value = "  -1.84%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +0.43%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.576.75"

# Row 16
$ws.Cells.Item(16, 4).Value = "63.268.69"
$ws.Cells.Item(16, 5).Value = "  -1.30%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +0.21%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "3.083.91"
$ws.Cells.Item(18, 5).Value = "  +0.43%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "510.45"
$ws.Cells.Item(19, 5).Value = "  +4.62%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "6.71"
$ws.Cells.Item(20, 5).Value = "  +1.36%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "13.70"
$ws.Cells.Item(21, 5).Value = "  +0.77%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "0.711"
$ws.Cells.Item(22, 5).Value = "  +3.82%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "7.33"
$ws.Cells.Item(23, 5).Value = "  +1.12%  "

# Row 24
$ws.Cells.Item(24, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(24, 4).Value = "12.44"
$ws.Cells.Item(24, 5).Value = "  +0.11%  "

# Row 25
$ws.Cells.Item(25, 2).Value = "Litecoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(25, 4).Value = "78.07"
$ws.Cells.Item(25, 5).Value = "  +0.13%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.02%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "2.78"
$ws.Cells.Item(27, 5).Value = "  +2.43%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "8.35"
$ws.Cells.Item(28, 5).Value = "  +0.73%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "2.05"
$ws.Cells.Item(29, 5).Value = "  -1.31%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "1.00"
$ws.Cells.Item(30, 5).Value = "  +0.02%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "26.37"
$ws.Cells.Item(31, 5).Value = "  +2.33%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "2.53"
$ws.Cells.Item(32, 5).Value = "  -4.22%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "1.12"
$ws.Cells.Item(33, 5).Value = "  -1.94%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "537.00"
$ws.Cells.Item(34, 5).Value = "  -9.14%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "58.05"
$ws.Cells.Item(35, 5).Value = "  +11.40%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "5.94"
$ws.Cells.Item(36, 5).Value = "  -0.28%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "5.18"
$ws.Cells.Item(37, 5).Value = "  -4.00%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "0.0416"
$ws.Cells.Item(38, 5).Value = "  +4.26%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "0.0800"
$ws.Cells.Item(39, 5).Value = "  +1.26%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "3.092.29"
$ws.Cells.Item(40, 5).Value = "  +4.07%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "0.119"
$ws.Cells.Item(41, 5).Value = "  +0.14%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "8.14"
$ws.Cells.Item(42, 5).Value = "  -0.77%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "2.67"
$ws.Cells.Item(43, 5).Value = "  -6.87%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "0.256"
$ws.Cells.Item(44, 5).Value = "  +4.65%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +0.04%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "2.10"
$ws.Cells.Item(46, 5).Value = "  +0.51%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "121.42"
$ws.Cells.Item(47, 5).Value = "  +1.68%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "24.16"
$ws.Cells.Item(48, 5).Value = "  -3.16%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +0.05%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "0.0₃0502"
$ws.Cells.Item(50, 5).Value = "  -5.42%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "2.34"
$ws.Cells.Item(51, 5).Value = "  +59.37%  "
